# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.841.21'
$ws.Range("E2").Value = '  +2.85%  '

$ws.Range("D3").Value = '3.200.63'
$ws.Range("E3").Value = '  +1.69%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''605.03'
$ws.Range("E5").Value = '  +4.49%  '

$ws.Range("D6").Value = '''157.09'
$ws.Range("E6").Value = '  +5.58%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '''0.556'
$ws.Range("E8").Value = '  +5.96%  '

$ws.Range("D9").Value = '3.198.16'
$ws.Range("E9").Value = '  +1.56%  '

$ws.Range("E10").Value = '  +1.88%  '

$ws.Range("D11").Value = '''5.92'
$ws.Range("E11").Value = '  -3.56%  '

$ws.Range("D12").Value = '''0.516'
$ws.Range("E12").Value = '  +3.40%  '

$ws.Range("E13").Value = '  +2.27%  '

$ws.Range("D14").Value = '''39.20'
$ws.Range("E14").Value = '  +5.68%  '

$ws.Range("D15").Value = '3.723.64'
$ws.Range("E15").Value = '  +1.62%  '

$ws.Range("D16").Value = '66.754.37'
$ws.Range("E16").Value = '  +2.81%  '

$ws.Range("D17").Value = '''7.47'
$ws.Range("E17").Value = '  +4.85%  '

$ws.Range("D18").Value = '3.201.24'
$ws.Range("E18").Value = '  +1.95%  '

$ws.Range("E19").Value = '  +0.87%  '

$ws.Range("D20").Value = '''522.82'
$ws.Range("E20").Value = '  +4.07%  '

$ws.Range("D21").Value = '''15.47'
$ws.Range("E21").Value = '  +2.69%  '

$ws.Range("D22").Value = '''0.743'
$ws.Range("E22").Value = '  +4.13%  '

$ws.Range("D23").Value = '''8.21'
$ws.Range("E23").Value = '  +6.37%  '

$ws.Range("D24").Value = '''15.09'
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = '''85.40'
$ws.Range("E25").Value = '  +1.44%  '

$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("D27").Value = '''9.27'
$ws.Range("E27").Value = '  +2.39%  '

$ws.Range("E28").Value = '  +3.89%  '

$ws.Range("D29").Value = '''2.39'
$ws.Range("E29").Value = '  +9.79%  '

$ws.Range("D30").Value = '''3.01'
$ws.Range("E30").Value = '  +7.66%  '

$ws.Range("D31").Value = '''7.02'
$ws.Range("E31").Value = '  +9.08%  '

$ws.Range("D32").Value = '''28.37'
$ws.Range("E32").Value = '  +3.29%  '

$ws.Range("E33").Value = '  +3.44%  '

$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("E35").Value = '  +1.61%  '

$ws.Range("D36").Value = '''519.83'
$ws.Range("E36").Value = '  +9.47%  '

$ws.Range("D37").Value = '''55.16'
$ws.Range("E37").Value = '  +0.36%  '

$ws.Range("D38").Value = '''0.0904'
$ws.Range("E38").Value = '  +2.00%  '

$ws.Range("D39").Value = '''0.0427'
$ws.Range("E39").Value = '  +3.50%  '

$ws.Range("E40").Value = '  +8.60%  '

$ws.Range("E41").Value = '  +2.25%  '

$ws.Range("E42").Value = '  -0.66%  '

$ws.Range("D43").Value = '0.0₃0689'
$ws.Range("E43").Value = '  +15.63%  '

$ws.Range("E44").Value = '  +7.37%  '

$ws.Range("D45").Value = '''2.47'
$ws.Range("E45").Value = '  +2.31%  '

$ws.Range("D46").Value = '2.896.93'
$ws.Range("E46").Value = '  -3.44%  '

$ws.Range("D47").Value = '''28.66'
$ws.Range("E47").Value = '  +1.52%  '

$ws.Range("D48").Value = '''2.75'
$ws.Range("E48").Value = '  +10.98%  '

$ws.Range("E49").Value = '  +3.85%  '

$ws.Range("E51").Value = '  +4.36%  '
